$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 64
$ws.Range("H64").Value = 3250
$ws.Range("J64").Value = 3500
$ws.Range("L64").Value = 3500
$ws.Range("N64").Value = -3996
# Row 67
$ws.Range("H67").Value = 3250
$ws.Range("J67").Value = 3500
$ws.Range("L67").Value = 3500
$ws.Range("N67").Value = -5216
# Row 88
$ws.Range("H88").Value = 5475
$ws.Range("J88").Value = 5475
$ws.Range("L88").Value = 5475
$ws.Range("N88").Value = -6287
# Row 91
$ws.Range("H91").Value = 5475
$ws.Range("J91").Value = 5475
$ws.Range("L91").Value = 5475
$ws.Range("N91").Value = -8283
# Row 106
$ws.Range("H106").Value = 3214.2144
$ws.Range("I106").Value = 2165.5833
$ws.Range("K106").Value = 2165.5833
$ws.Range("M106").Value = -1534.5833
# Row 107
$ws.Range("H107").Value = 846
$ws.Range("I107").Value = 680.7143
$ws.Range("K107").Value = 680.7143
$ws.Range("M107").Value = 1239.2857
# Row 116
$ws.Range("H116").Value = 8299
$ws.Range("I116").Value = 11279.7
$ws.Range("K116").Value = 11279.7
$ws.Range("M116").Value = -7837.700000000001
# Row 132
$ws.Range("H132").Value = 1184.6857
$ws.Range("I132").Value = 998.875
$ws.Range("J132").Value = 3166.6667
$ws.Range("K132").Value = 2996.625
$ws.Range("L132").Value = 9500.000100000001
$ws.Range("M132").Value = -466.625
$ws.Range("N132").Value = -14560.0001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3289.7454
$ws.Range("I32").Value = 2088.3901
$ws.Range("K32").Value = 2088.3901
$ws.Range("M32").Value = -1801.3901
# Row 45
$ws.Range("H45").Value = 1820.6
$ws.Range("I45").Value = 1768.3334
$ws.Range("K45").Value = 1768.3334
$ws.Range("M45").Value = -1391.3334
# Row 61
$ws.Range("H61").Value = 3984.8096
$ws.Range("I61").Value = 3030.75
$ws.Range("J61").Value = 7037.8
$ws.Range("K61").Value = 3030.75
$ws.Range("L61").Value = 7037.8
$ws.Range("M61").Value = -2818.75
$ws.Range("N61").Value = -7461.8
# Row 74
$ws.Range("H74").Value = 3518.0667
$ws.Range("I74").Value = 3443.8462
$ws.Range("K74").Value = 3443.8462
$ws.Range("M74").Value = -2569.8462
# Row 77
$ws.Range("H77").Value = 3518.0667
$ws.Range("I77").Value = 3443.8462
$ws.Range("K77").Value = 17219.231
$ws.Range("M77").Value = -12851.231
# Row 136
$ws.Range("H136").Value = 3984.8096
$ws.Range("I136").Value = 3030.75
$ws.Range("J136").Value = 7037.8
$ws.Range("K136").Value = 9092.25
$ws.Range("L136").Value = 21113.4
$ws.Range("M136").Value = -6542.25
$ws.Range("N136").Value = -26213.4

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 107182.58
$ws.Range("I86").Value = 2041.2142
$ws.Range("J86").Value = 401578.4
$ws.Range("K86").Value = 2041.2142
$ws.Range("L86").Value = 401578.4
$ws.Range("M86").Value = -918.2141999999999
$ws.Range("N86").Value = -403824.4
# Row 89
$ws.Range("H89").Value = 107182.58
$ws.Range("I89").Value = 2041.2142
$ws.Range("J89").Value = 401578.4
$ws.Range("K89").Value = 10206.071
$ws.Range("L89").Value = 2007892
$ws.Range("M89").Value = -4590.071
$ws.Range("N89").Value = -2019124
# Row 99
$ws.Range("H99").Value = 1570.4286
$ws.Range("I99").Value = 1498.5
$ws.Range("J99").Value = 1666.3334
$ws.Range("K99").Value = 1498.5
$ws.Range("L99").Value = 1666.3334
$ws.Range("M99").Value = -0.5
$ws.Range("N99").Value = -4662.3334
# Row 105
$ws.Range("H105").Value = 2272.2104
$ws.Range("I105").Value = 2345.5334
$ws.Range("K105").Value = 2345.5334
$ws.Range("M105").Value = -598.5333999999998
# Row 107
$ws.Range("H107").Value = 1814.875
$ws.Range("I107").Value = 1468
$ws.Range("J107").Value = 3318
$ws.Range("K107").Value = 1468
$ws.Range("L107").Value = 3318
$ws.Range("M107").Value = 452
$ws.Range("N107").Value = -7158

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1000
$ws.Range("I16").Value = 1000
$ws.Range("K16").Value = 1000
$ws.Range("M16").Value = -713
# Row 31
$ws.Range("H31").Value = 3062.8333
$ws.Range("I31").Value = 900.5333000000001
$ws.Range("J31").Value = 6666.6665
$ws.Range("K31").Value = 900.5333000000001
$ws.Range("L31").Value = 6666.6665
$ws.Range("M31").Value = -605.5333000000001
$ws.Range("N31").Value = -7256.6665
# Row 34
$ws.Range("H34").Value = 3062.8333
$ws.Range("I34").Value = 900.5333000000001
$ws.Range("J34").Value = 6666.6665
$ws.Range("K34").Value = 900.5333000000001
$ws.Range("L34").Value = 6666.6665
$ws.Range("M34").Value = -698.5333000000001
$ws.Range("N34").Value = -7070.6665
# Row 52
$ws.Range("H52").Value = 64280
$ws.Range("J52").Value = 64280
$ws.Range("L52").Value = 64280
$ws.Range("N52").Value = -64868
# Row 113
$ws.Range("H113").Value = 1000
$ws.Range("I113").Value = 1000
$ws.Range("K113").Value = 1000
$ws.Range("M113").Value = 1170

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 7928.31
$ws.Range("I131").Value = 667.5
$ws.Range("J131").Value = 8391.766
$ws.Range("K131").Value = 2002.5
$ws.Range("L131").Value = 25175.298
$ws.Range("M131").Value = 3037.5
$ws.Range("N131").Value = -35255.298

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 10
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").ClearContents()
# Row 14
$ws.Range("H14").Value = 6900000
$ws.Range("I14").Value = 6900000
$ws.Range("K14").Value = 6900000
$ws.Range("M14").Value = -6899832
# Row 70
$ws.Range("H70").Value = 4442.6665
$ws.Range("I70").Value = 4608.625
$ws.Range("J70").Value = 4110.75
$ws.Range("K70").Value = 4608.625
$ws.Range("L70").Value = 4110.75
$ws.Range("M70").Value = -4338.625
$ws.Range("N70").Value = -4650.75
# Row 73
$ws.Range("H73").Value = 4442.6665
$ws.Range("I73").Value = 4608.625
$ws.Range("J73").Value = 4110.75
$ws.Range("K73").Value = 4608.625
$ws.Range("L73").Value = 4110.75
$ws.Range("M73").Value = -3672.625
$ws.Range("N73").Value = -5982.75

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 55
$ws.Range("H55").Value = 285.16666
$ws.Range("I55").Value = 151.25
$ws.Range("K55").Value = 151.25
$ws.Range("M55").Value = 21.75
# Row 100
$ws.Range("H100").Value = 2305.5
$ws.Range("J100").Value = 3000
$ws.Range("L100").Value = 3000
$ws.Range("N100").Value = -4082

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 51
$ws.Range("H51").Value = 34344
$ws.Range("J51").Value = 34344
$ws.Range("L51").Value = 34344
$ws.Range("N51").Value = -35364
# Row 107
$ws.Range("H107").Value = 558.125
$ws.Range("I107").Value = 416.26315
$ws.Range("K107").Value = 1248.78945
$ws.Range("M107").Value = 671.21055
# Row 123
$ws.Range("H123").Value = 47518.9
$ws.Range("J123").Value = 47518.9
$ws.Range("L123").Value = 47518.9
$ws.Range("N123").Value = -57318.9
# Row 135
$ws.Range("H135").Value = 95335.2
$ws.Range("J135").Value = 95335.2
$ws.Range("L135").Value = 95335.2
$ws.Range("N135").Value = -105475.2
